$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# EFECTO LE S.A.S. (NIT 9013604463) - estado de cuenta
# Elimina periodos de mora anteriores y agrega los nuevos periodos,
# reordenando cada trabajador con su periodo mas reciente primero.

$data = @(
    @("CC", "1143386199", "DAYANA CAROLINA PEREZ HURTADO", "2408", 52000, 1300000),
    @("CC", "1143386199", "DAYANA CAROLINA PEREZ HURTADO", "2407", 52000, 1300000),
    @("CC", "1143386199", "DAYANA CAROLINA PEREZ HURTADO", "2406", 52000, 1300000),
    @("CC", "1143386199", "DAYANA CAROLINA PEREZ HURTADO", "2405", 52000, 1300000),
    @("CC", "1143386199", "DAYANA CAROLINA PEREZ HURTADO", "2404", 52000, 1300000),
    @("CC", "92070402", "ORLANDO RAMON BELTRAN RODRIGUEZ", "2410", 15600, 1300000),
    @("CC", "92070402", "ORLANDO RAMON BELTRAN RODRIGUEZ", "2409", 52000, 1300000),
    @("CC", "92070402", "ORLANDO RAMON BELTRAN RODRIGUEZ", "2408", 52000, 1300000),
    @("CC", "92070402", "ORLANDO RAMON BELTRAN RODRIGUEZ", "2407", 52000, 1300000),
    @("CC", "92070402", "ORLANDO RAMON BELTRAN RODRIGUEZ", "2406", 52000, 1300000),
    @("CC", "92070402", "ORLANDO RAMON BELTRAN RODRIGUEZ", "2405", 52000, 1300000),
    @("CC", "92070402", "ORLANDO RAMON BELTRAN RODRIGUEZ", "2404", 52000, 1300000)
)

$row = 16
foreach ($r in $data) {
    $ws.Cells.Item($row, 2).Value = $r[0]
    $ws.Cells.Item($row, 3).Value = $r[1]
    $ws.Cells.Item($row, 4).Value = $r[2]
    $ws.Cells.Item($row, 5).Value = $r[3]
    $ws.Cells.Item($row, 6).Value = $r[4]
    $ws.Cells.Item($row, 7).Value = $r[5]
    $row = $row + 1
}
